# issue #5: stock data output to json file
# Add a "property_category" column to the 股票 (stock) sheet, populated
# with the literal value "stock" for the existing data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("股票")

# Insert a new column before the existing "date" column (H) so that the
# subsequent columns (date, legislator_name, legislator_id) shift right
# by one, matching the target layout:
#   B name | C owner | D quantity | E face_value | F currency | G total |
#   H property_category | I date | J legislator_name | K legislator_id
$ws.Columns.Item(8).Insert()

$ws.Cells.Item(1, 8).Value = "property_category"
$ws.Cells.Item(2, 8).Value = "stock"
